# Simplify runs: merge adjacent text runs that share identical formatting
# into a single run (pure "re-split" edit - the rendered text itself does
# not change, only how it is broken up across <a:r> elements).
#
# We do this with TextRange.Characters(Start, Length) sub-ranges: selecting
# the exact span that covers the two runs-to-merge and re-assigning the very
# same text to it collapses that span into a single run using the formatting
# of the first run in the span, while runs outside the span are left
# untouched.

$p = $ppt.ActivePresentation

# --- Slide 8: title "前序遍历:（中左右）: 6-4-2-1-0-3-8-7-9-10<br>中序遍历:
#     （左中右）: 0-1-2-3-4-6-7-8-9-10<br>后序遍历:（左右中）: 0-1-3-2-4-7-10-9-8-6<br>"
$slide8 = $p.Slides.Item(8)
$title8 = $slide8.Shapes.Item(1).TextFrame.TextRange

# Merge ": " + "0-1-2-3-4-6-7-8-9-10" -> ": 0-1-2-3-4-6-7-8-9-10"
$title8.Characters(44, 22).Text = ": 0-1-2-3-4-6-7-8-9-10"

# Merge "（左右中" + "）" -> "（左右中）"
$title8.Characters(72, 5).Text = "（左右中）"

# --- Slide 9: title "非递归前序遍历: 6-4-2-1-0-3-8-7-9-10"
$slide9 = $p.Slides.Item(9)
$title9 = $slide9.Shapes.Item(1).TextFrame.TextRange

# Merge "非递归" + "前序遍历" -> "非递归前序遍历"
$title9.Characters(1, 7).Text = "非递归前序遍历"

# --- Slide 10: title "后序遍历: 0-1-3-2-4-7-10-9-8-6"
$slide10 = $p.Slides.Item(10)
$title10 = $slide10.Shapes.Item(1).TextFrame.TextRange

# Merge "后" + "序遍历" -> "后序遍历"
$title10.Characters(1, 4).Text = "后序遍历"
